$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 16
$ws1.Range("F7").Value = 3861
$ws1.Range("F8").Value = 2551
$ws1.Range("F10").Value = 2416
$ws1.Range("F14").Value = 1630
$ws1.Range("F15").Value = 645
$ws1.Range("F16").Value = 8
$ws1.Range("F17").Value = 97
$ws1.Range("F20").Value = 51
$ws1.Range("F22").Value = 67
$ws1.Range("F23").Value = 442
$ws1.Range("F27").Value = 672
$ws1.Range("F28").Value = 86
$ws1.Range("F29").Value = 74
$ws1.Range("F30").Value = 373
$ws1.Range("F31").Value = 39
$ws1.Range("F33").Value = 881
$ws1.Range("F34").Value = 46
$ws1.Range("F35").Value = 9
$ws1.Range("F36").Value = 941
$ws1.Range("F37").Value = 1959
$ws1.Range("F38").Value = 227
$ws1.Range("F39").Value = 519
$ws1.Range("F41").Value = 9
$ws1.Range("F42").Value = 598
$ws1.Range("F43").Value = 1250
$ws1.Range("F44").Value = 42
$ws1.Range("F46").Value = 415

# --- Sheet: 演出 ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 64

# --- Sheet: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 16
$ws4.Range("F6").Value = 3861
$ws4.Range("F7").Value = 2551
$ws4.Range("F8").Value = 2416
$ws4.Range("F10").Value = 1630
$ws4.Range("F12").Value = 645
$ws4.Range("F13").Value = 8
$ws4.Range("F14").Value = 97
$ws4.Range("F17").Value = 51
$ws4.Range("F19").Value = 67
$ws4.Range("F20").Value = 442
$ws4.Range("F24").Value = 672
$ws4.Range("F25").Value = 86
$ws4.Range("F26").Value = 64
$ws4.Range("F29").Value = 74
$ws4.Range("F30").Value = 373
$ws4.Range("F32").Value = 881
$ws4.Range("F33").Value = 46
$ws4.Range("F34").Value = 9
$ws4.Range("F36").Value = 942
$ws4.Range("F37").Value = 1959
$ws4.Range("F38").Value = 227
$ws4.Range("F42").Value = 519
$ws4.Range("F44").Value = 9
$ws4.Range("F45").Value = 598
$ws4.Range("F46").Value = 1250
$ws4.Range("F47").Value = 42
$ws4.Range("F48").Value = 415
